$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data_entities")
$ws2 = $wb.Worksheets.Item("functional_diagram")
$ws1.Range("A1:E30").Copy()
$ws2.Range("A15").PasteSpecial(-4163)
$ws1.Range("A1:E30").Copy()
$ws2.Range("A15").PasteSpecial(-4122)
$ws2.Columns.Item(1).AutoFit() | Out-Null
$ws2.Columns.Item(2).AutoFit() | Out-Null
$ws2.Columns.Item(4).AutoFit() | Out-Null
$ws2.Columns.Item(5).AutoFit() | Out-Null
Write-Host "ColA width: $($ws2.Columns.Item(1).ColumnWidth)"
Write-Host "ColB width: $($ws2.Columns.Item(2).ColumnWidth)"
Write-Host "ColD width: $($ws2.Columns.Item(4).ColumnWidth)"
Write-Host "ColE width: $($ws2.Columns.Item(5).ColumnWidth)"
